{"js": "// Apply the dated worksheet refresh: update the date line and all\n// \"dividend\u00f7divisor=quotient, remainder\" answer cells in the table.\nconst replacements = [\n  [\"2024-10-12 Saturday\", \"2024-10-13 Sunday\"],\n  [\"994\u00f72=497, 0\", \"898\u00f77=128, 2\"],\n  [\"242\u00f75=48, 2\", \"674\u00f75=134, 4\"],\n  [\"645\u00f78=80, 5\", \"377\u00f75=75, 2\"],\n  [\"767\u00f78=95, 7\", \"412\u00f72=206, 0\"],\n  [\"679\u00f79=75, 4\", \"920\u00f78=115, 0\"],\n  [\"555\u00f77=79, 2\", \"332\u00f79=36, 8\"],\n  [\"432\u00f75=86, 2\", \"360\u00f78=45, 0\"],\n  [\"611\u00f72=305, 1\", \"277\u00f78=34, 5\"],\n  [\"240\u00f78=30, 0\", \"376\u00f72=188, 0\"],\n  [\"956\u00f77=136, 4\", \"198\u00f79=22, 0\"],\n  [\"237\u00f79=26, 3\", \"887\u00f75=177, 2\"],\n  [\"356\u00f78=44, 4\", \"160\u00f76=26, 4\"],\n  [\"298\u00f74=74, 2\", \"562\u00f73=187, 1\"],\n  [\"781\u00f78=97, 5\", \"326\u00f76=54, 2\"],\n  [\"123\u00f74=30, 3\", \"127\u00f77=18, 1\"],\n  [\"267\u00f78=33, 3\", \"570\u00f76=95, 0\"],\n  [\"926\u00f78=115, 6\", \"651\u00f74=162, 3\"],\n  [\"365\u00f73=121, 2\", \"100\u00f72=50, 0\"],\n  [\"150\u00f72=75, 0\", \"289\u00f72=144, 1\"],\n  [\"255\u00f75=51, 0\", \"649\u00f76=108, 1\"],\n  [\"559\u00f79=62, 1\", \"817\u00f79=90, 7\"],\n  [\"248\u00f75=49, 3\", \"145\u00f72=72, 1\"],\n  [\"884\u00f76=147, 2\", \"652\u00f74=163, 0\"],\n  [\"730\u00f73=243, 1\", \"363\u00f77=51, 6\"],\n  [\"558\u00f76=93, 0\", \"652\u00f72=326, 0\"]\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the dated worksheet refresh: update the date line and all\n# \"dividend\u00f7divisor=quotient, remainder\" answer cells in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-12 Saturday\", \"2024-10-13 Sunday\"),\n    @(\"994\u00f72=497, 0\", \"898\u00f77=128, 2\"),\n    @(\"242\u00f75=48, 2\", \"674\u00f75=134, 4\"),\n    @(\"645\u00f78=80, 5\", \"377\u00f75=75, 2\"),\n    @(\"767\u00f78=95, 7\", \"412\u00f72=206, 0\"),\n    @(\"679\u00f79=75, 4\", \"920\u00f78=115, 0\"),\n    @(\"555\u00f77=79, 2\", \"332\u00f79=36, 8\"),\n    @(\"432\u00f75=86, 2\", \"360\u00f78=45, 0\"),\n    @(\"611\u00f72=305, 1\", \"277\u00f78=34, 5\"),\n    @(\"240\u00f78=30, 0\", \"376\u00f72=188, 0\"),\n    @(\"956\u00f77=136, 4\", \"198\u00f79=22, 0\"),\n    @(\"237\u00f79=26, 3\", \"887\u00f75=177, 2\"),\n    @(\"356\u00f78=44, 4\", \"160\u00f76=26, 4\"),\n    @(\"298\u00f74=74, 2\", \"562\u00f73=187, 1\"),\n    @(\"781\u00f78=97, 5\", \"326\u00f76=54, 2\"),\n    @(\"123\u00f74=30, 3\", \"127\u00f77=18, 1\"),\n    @(\"267\u00f78=33, 3\", \"570\u00f76=95, 0\"),\n    @(\"926\u00f78=115, 6\", \"651\u00f74=162, 3\"),\n    @(\"365\u00f73=121, 2\", \"100\u00f72=50, 0\"),\n    @(\"150\u00f72=75, 0\", \"289\u00f72=144, 1\"),\n    @(\"255\u00f75=51, 0\", \"649\u00f76=108, 1\"),\n    @(\"559\u00f79=62, 1\", \"817\u00f79=90, 7\"),\n    @(\"248\u00f75=49, 3\", \"145\u00f72=72, 1\"),\n    @(\"884\u00f76=147, 2\", \"652\u00f74=163, 0\"),\n    @(\"730\u00f73=243, 1\", \"363\u00f77=51, 6\"),\n    @(\"558\u00f76=93, 0\", \"652\u00f72=326, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    [void]$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
